$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.931.42'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.555.37'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '1.776.42'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '1.555.13'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').Value = '26.919.10'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.72'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '218.38'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('E30').Value = '  +2.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.09'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('D33').Value = '1.433.05'
$ws.Range('E33').Value = '  +4.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.06'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.34%  '
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.982'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.30%  '
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.519'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.27'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.51%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.985'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.86'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').Value = '1.690.59'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.90'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.15%  '
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('E50').Value = '  +3.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0953'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.29%  '
